# Add a second worksheet "헬스비급여_2" that duplicates the existing
# "헬스비급여" sheet (same header row + all data rows), matching the
# workbook's new second sheet (sheetId=2) added in this commit
# ("add more function for to_excel()").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Copy the sheet, placing the new copy immediately after the original.
$ws1.Copy($null, $ws1)

# Excel auto-names the copy "헬스비급여 (2)"; rename to match target.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "헬스비급여_2"
